$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 15-22, years 2013-2020 under old layout)
$ws.Range("A15:G22").EntireRow.Delete()

# Overwrite rows 2-14 with the refreshed 2010-2022 series
$ws.Cells.Item(2, 1).Value2 = "2010年"
$ws.Cells.Item(2, 2).Value2 = 30807.9326754466
$ws.Cells.Item(2, 3).Value2 = 412119.255796083
$ws.Cells.Item(2, 4).Value2 = 410354.11014066
$ws.Cells.Item(2, 5).Value2 = 38430.8498088599
$ws.Cells.Item(2, 6).Value2 = 182061.890339732
$ws.Cells.Item(2, 7).Value2 = 191626.515647491

$ws.Cells.Item(3, 1).Value2 = "2011年"
$ws.Cells.Item(3, 2).Value2 = 36277.1363217578
$ws.Cells.Item(3, 3).Value2 = 487940.180525355
$ws.Cells.Item(3, 4).Value2 = 483392.794744631
$ws.Cells.Item(3, 5).Value2 = 44781.4604971136
$ws.Cells.Item(3, 6).Value2 = 216123.620990973
$ws.Cells.Item(3, 7).Value2 = 227035.099037269

$ws.Cells.Item(4, 1).Value2 = "2012年"
$ws.Cells.Item(4, 2).Value2 = 39771.3728109773
$ws.Cells.Item(4, 3).Value2 = 538579.953468974
$ws.Cells.Item(4, 4).Value2 = 537329.007798286
$ws.Cells.Item(4, 5).Value2 = 49084.6359218683
$ws.Cells.Item(4, 6).Value2 = 244856.249002227
$ws.Cells.Item(4, 7).Value2 = 244639.068544879

$ws.Cells.Item(5, 1).Value2 = "2013年"
$ws.Cells.Item(5, 2).Value2 = 43496.6131824886
$ws.Cells.Item(5, 3).Value2 = 592963.229548957
$ws.Cells.Item(5, 4).Value2 = 588141.211190155
$ws.Cells.Item(5, 5).Value2 = 53028.0728537458
$ws.Cells.Item(5, 6).Value2 = 277983.542768454
$ws.Cells.Item(5, 7).Value2 = 261951.613926757

$ws.Cells.Item(6, 1).Value2 = "2014年"
$ws.Cells.Item(6, 2).Value2 = 46911.718728133
$ws.Cells.Item(6, 3).Value2 = 643563.104543766
$ws.Cells.Item(6, 4).Value2 = 644380.151914331
$ws.Cells.Item(6, 5).Value2 = 55626.3215706793
$ws.Cells.Item(6, 6).Value2 = 310653.963237863
$ws.Cells.Item(6, 7).Value2 = 277282.819735223

$ws.Cells.Item(7, 1).Value2 = "2015年"
$ws.Cells.Item(7, 2).Value2 = 49922.3267613586
$ws.Cells.Item(7, 3).Value2 = 688858.218049283
$ws.Cells.Item(7, 4).Value2 = 685571.21848939
$ws.Cells.Item(7, 5).Value2 = 57774.6412645268
$ws.Cells.Item(7, 6).Value2 = 349744.650043795
$ws.Cells.Item(7, 7).Value2 = 281338.926740961

$ws.Cells.Item(8, 1).Value2 = "2016年"
$ws.Cells.Item(8, 2).Value2 = 53782.9973903485
$ws.Cells.Item(8, 3).Value2 = 746395.059483517
$ws.Cells.Item(8, 4).Value2 = 742694.053654712
$ws.Cells.Item(8, 5).Value2 = 60139.196087384
$ws.Cells.Item(8, 6).Value2 = 390828.059989709
$ws.Cells.Item(8, 7).Value2 = 295427.803406424

$ws.Cells.Item(9, 1).Value2 = "2017年"
$ws.Cells.Item(9, 2).Value2 = 59592.2510902632
$ws.Cells.Item(9, 3).Value2 = 832035.948559918
$ws.Cells.Item(9, 4).Value2 = 830945.699737374
$ws.Cells.Item(9, 5).Value2 = 62099.5439820457
$ws.Cells.Item(9, 6).Value2 = 438355.94736846
$ws.Cells.Item(9, 7).Value2 = 331580.457209411

$ws.Cells.Item(10, 1).Value2 = "2018年"
$ws.Cells.Item(10, 2).Value2 = 65533.7426977278
$ws.Cells.Item(10, 3).Value2 = 919281.129066646
$ws.Cells.Item(10, 4).Value2 = 915243.45391614
$ws.Cells.Item(10, 5).Value2 = 64745.1561040139
$ws.Cells.Item(10, 6).Value2 = 489700.762579774
$ws.Cells.Item(10, 7).Value2 = 364835.210382857

$ws.Cells.Item(11, 1).Value2 = "2019年"
$ws.Cells.Item(11, 2).Value2 = 70077.6917901966
$ws.Cells.Item(11, 3).Value2 = 986515.202291903
$ws.Cells.Item(11, 4).Value2 = 983751.20256033
$ws.Cells.Item(11, 5).Value2 = 70473.5905931704
$ws.Cells.Item(11, 6).Value2 = 535370.99070492
$ws.Cells.Item(11, 7).Value2 = 380670.620993813

$ws.Cells.Item(12, 1).Value2 = "2020年"
$ws.Cells.Item(12, 2).Value2 = 71828.1484112169
$ws.Cells.Item(12, 3).Value2 = 1013567.00223068
$ws.Cells.Item(12, 4).Value2 = 1005451.31315167
$ws.Cells.Item(12, 5).Value2 = 78030.9003457103
$ws.Cells.Item(12, 6).Value2 = 551973.748218789
$ws.Cells.Item(12, 7).Value2 = 383562.353666183

$ws.Cells.Item(13, 1).Value2 = "2021年"
$ws.Cells.Item(13, 2).Value2 = 81370
$ws.Cells.Item(13, 3).Value2 = 1149237
$ws.Cells.Item(13, 4).Value2 = 1141230.8
$ws.Cells.Item(13, 5).Value2 = 83216.5
$ws.Cells.Item(13, 6).Value2 = 614476.4
$ws.Cells.Item(13, 7).Value2 = 451544.1

$ws.Cells.Item(14, 1).Value2 = "2022年"
$ws.Cells.Item(14, 2).Value2 = 85698
$ws.Cells.Item(14, 3).Value2 = 1210207.2
$ws.Cells.Item(14, 4).Value2 = 1197250.4
$ws.Cells.Item(14, 5).Value2 = 88345.1
$ws.Cells.Item(14, 6).Value2 = 638697.6
$ws.Cells.Item(14, 7).Value2 = 483164.5

